# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the Hyperion_Profits workbook (columns H:N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2262.5
$ws.Range("I2").Value = 2699.8
$ws.Range("J2").Value = 1533.6666
$ws.Range("K2").Value = 2699.8
$ws.Range("L2").Value = 1533.6666
$ws.Range("M2").Value = -2586.8
$ws.Range("N2").Value = -1759.6666

$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2825
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 5499.909
$ws.Range("J64").Value = 8375
$ws.Range("L64").Value = 8375
$ws.Range("N64").Value = -8871

$ws.Range("H67").Value = 5499.909
$ws.Range("J67").Value = 8375
$ws.Range("L67").Value = 8375
$ws.Range("N67").Value = -10091

$ws.Range("H137").Value = 66880.78999999999
$ws.Range("I137").Value = 114090.125
$ws.Range("K137").Value = 342270.375
$ws.Range("M137").Value = -339720.375

$ws.Range("H141").Value = 1984.75
$ws.Range("I141").Value = 1984.75
$ws.Range("K141").Value = 5954.25
$ws.Range("M141").Value = -774.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1614.8823
$ws.Range("I2").Value = 1049.1666
$ws.Range("J2").Value = 2972.6
$ws.Range("K2").Value = 1049.1666
$ws.Range("L2").Value = 2972.6
$ws.Range("M2").Value = -936.1666
$ws.Range("N2").Value = -3198.6

$ws.Range("H61").Value = 1921.8788
$ws.Range("I61").Value = 1424.3529
$ws.Range("K61").Value = 1424.3529
$ws.Range("M61").Value = -1212.3529

$ws.Range("H74").Value = 47403.348
$ws.Range("I74").Value = 5946.4443
$ws.Range("J74").Value = 260610.28
$ws.Range("K74").Value = 5946.4443
$ws.Range("L74").Value = 260610.28
$ws.Range("M74").Value = -5072.4443
$ws.Range("N74").Value = -262358.28

$ws.Range("H77").Value = 47403.348
$ws.Range("I77").Value = 5946.4443
$ws.Range("J77").Value = 260610.28
$ws.Range("K77").Value = 29732.2215
$ws.Range("L77").Value = 1303051.4
$ws.Range("M77").Value = -25364.2215
$ws.Range("N77").Value = -1311787.4

$ws.Range("H116").Value = 1614.8823
$ws.Range("I116").Value = 1049.1666
$ws.Range("J116").Value = 2972.6
$ws.Range("K116").Value = 1049.1666
$ws.Range("L116").Value = 2972.6
$ws.Range("M116").Value = 1244.8334
$ws.Range("N116").Value = -7560.6

$ws.Range("H136").Value = 1921.8788
$ws.Range("I136").Value = 1424.3529
$ws.Range("K136").Value = 4273.0587
$ws.Range("M136").Value = -1723.0587


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1614.8823
$ws.Range("I3").Value = 1049.1666
$ws.Range("J3").Value = 2972.6
$ws.Range("K3").Value = 1049.1666
$ws.Range("L3").Value = 2972.6
$ws.Range("M3").Value = -935.1666
$ws.Range("N3").Value = -3200.6


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 3077442.5
$ws.Range("I19").Value = 3333729.5
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 3333729.5
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -3333559.5
$ws.Range("N19").Value = -2340

$ws.Range("H23").Value = 26875
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4760

$ws.Range("H24").Value = 3077442.5
$ws.Range("I24").Value = 3333729.5
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 3333729.5
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = -3333559.5
$ws.Range("N24").Value = -2340

$ws.Range("H27").Value = 26875
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4808

$ws.Range("H31").Value = 30427.688
$ws.Range("I31").Value = 1425.2354
$ws.Range("K31").Value = 1425.2354
$ws.Range("M31").Value = -1130.2354

$ws.Range("H34").Value = 30427.688
$ws.Range("I34").Value = 1425.2354
$ws.Range("K34").Value = 1425.2354
$ws.Range("M34").Value = -1223.2354

$ws.Range("H58").Value = 1858.7826
$ws.Range("I58").Value = 1524.375
$ws.Range("K58").Value = 1524.375
$ws.Range("M58").Value = -1321.375

$ws.Range("H109").Value = 27371.5
$ws.Range("J109").Value = 27371.5
$ws.Range("L109").Value = 27371.5
$ws.Range("N109").Value = -29451.5

$ws.Range("H132").Value = 51993.61
$ws.Range("I132").Value = 2592.6428
$ws.Range("J132").Value = 224897
$ws.Range("K132").Value = 7777.928400000001
$ws.Range("L132").Value = 674691
$ws.Range("M132").Value = -5247.928400000001
$ws.Range("N132").Value = -679751

$ws.Range("H134").Value = 3784.3572
$ws.Range("I134").Value = 2721.5
$ws.Range("K134").Value = 8164.5
$ws.Range("M134").Value = -5629.5

$ws.Range("H136").Value = 1858.7826
$ws.Range("I136").Value = 1524.375
$ws.Range("K136").Value = 4573.125
$ws.Range("M136").Value = -2023.125


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5530369.5
$ws.Range("I4").Value = 6140494
$ws.Range("J4").Value = 446000
$ws.Range("K4").Value = 18421482
$ws.Range("L4").Value = 1338000
$ws.Range("M4").Value = -18421370
$ws.Range("N4").Value = -1338224

$ws.Range("H5").Value = 1112.2222
$ws.Range("I5").Value = 737
$ws.Range("J5").Value = 1658
$ws.Range("K5").Value = 2211
$ws.Range("L5").Value = 4974
$ws.Range("M5").Value = -2099
$ws.Range("N5").Value = -5198

$ws.Range("H103").Value = 149.5
$ws.Range("I103").Value = 99.666664
$ws.Range("K103").Value = 298.999992
$ws.Range("M103").Value = 580.000008

$ws.Range("H135").Value = 1112.2222
$ws.Range("I135").Value = 737
$ws.Range("J135").Value = 1658
$ws.Range("K135").Value = 6633
$ws.Range("L135").Value = 14922
$ws.Range("M135").Value = -4098
$ws.Range("N135").Value = -19992


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7298.077
$ws.Range("I70").Value = 7988.6
$ws.Range("K70").Value = 7988.6
$ws.Range("M70").Value = -7718.6

$ws.Range("H73").Value = 7298.077
$ws.Range("I73").Value = 7988.6
$ws.Range("K73").Value = 7988.6
$ws.Range("M73").Value = -7052.6


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 9504.5
$ws.Range("I4").Value = 9009
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 9009
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -8896
$ws.Range("N4").Value = -10226

$ws.Range("H28").Value = 9504.5
$ws.Range("I28").Value = 9009
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 9009
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -8777
$ws.Range("N28").Value = -10464

$ws.Range("H37").Value = 9504.5
$ws.Range("I37").Value = 9009
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 9009
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -8902
$ws.Range("N37").Value = -10214

$ws.Range("H68").Value = 2777.7144
$ws.Range("J68").Value = 5125
$ws.Range("L68").Value = 5125
$ws.Range("N68").Value = -6623

$ws.Range("H71").Value = 2777.7144
$ws.Range("J71").Value = 5125
$ws.Range("L71").Value = 25625
$ws.Range("N71").Value = -33113


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 804.7083
$ws.Range("I113").Value = 585.93335
$ws.Range("J113").Value = 1169.3334
$ws.Range("K113").Value = 1757.80005
$ws.Range("L113").Value = 3508.0002
$ws.Range("M113").Value = 412.1999499999999
$ws.Range("N113").Value = -7848.0002

$ws.Range("H132").Value = 55492.156
$ws.Range("I132").Value = 9145.177
$ws.Range("K132").Value = 27435.531
$ws.Range("M132").Value = -24905.531

